$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 332.22223
$ws.Range("I18").Value = 348.75
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 348.75
$ws.Range("L18").Value = 200
$ws.Range("M18").Value = -64.75
$ws.Range("N18").Value = -768

# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 1900
$ws.Range("I43").Value = 1166.6666
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 1166.6666
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = -1097.6666
$ws.Range("N43").Value = -3138

# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 1337.0714
$ws.Range("I62").Value = 1247.1818
$ws.Range("K62").Value = 1247.1818
$ws.Range("M62").Value = -623.1818000000001

# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 1337.0714
$ws.Range("I65").Value = 1247.1818
$ws.Range("K65").Value = 6235.909000000001
$ws.Range("M65").Value = -3115.909000000001

# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 3262.4546
$ws.Range("I86").Value = 1788.3
$ws.Range("K86").Value = 1788.3
$ws.Range("M86").Value = -665.3

# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 3262.4546
$ws.Range("I89").Value = 1788.3
$ws.Range("K89").Value = 8941.5
$ws.Range("M89").Value = -3325.5

# Row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 11785418
$ws.Range("I92").Value = 2416061.8
$ws.Range("J92").Value = 33334936
$ws.Range("K92").Value = 2416061.8
$ws.Range("L92").Value = 33334936
$ws.Range("M92").Value = -2414813.8
$ws.Range("N92").Value = -33337432

# Row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 16667610
$ws.Range("J112").Value = 23810568
$ws.Range("L112").Value = 71431704
$ws.Range("N112").Value = -71433920

# Row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 7863.7144
$ws.Range("I113").Value = 6000
$ws.Range("J113").Value = 8007.077
$ws.Range("K113").Value = 6000
$ws.Range("L113").Value = 8007.077
$ws.Range("M113").Value = -2746
$ws.Range("N113").Value = -14515.077

# Row 115 (Leve Item ID 27957)
$ws.Range("H115").Value = 580.8
$ws.Range("I115").Value = 580.8
$ws.Range("K115").Value = 1742.4
$ws.Range("M115").Value = -175.3999999999999

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 7052.8945
$ws.Range("I116").Value = 9883.75
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 9883.75
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = -6441.75
$ws.Range("N116").Value = -9084

# Row 125 (Leve Item ID 36228)
$ws.Range("H125").Value = 5424
$ws.Range("I125").Value = 8288
$ws.Range("J125").Value = 4350
$ws.Range("K125").Value = 74592
$ws.Range("L125").Value = 39150
$ws.Range("M125").Value = -72132
$ws.Range("N125").Value = -44070

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 1441.8889
$ws.Range("I132").Value = 1640.5652
$ws.Range("J132").Value = 299.5
$ws.Range("K132").Value = 4921.6956
$ws.Range("L132").Value = 898.5
$ws.Range("M132").Value = -2391.6956
$ws.Range("N132").Value = -5958.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 10 (Leve Item ID 2662)
$ws.Range("H10").Value = 70003.75
$ws.Range("J10").Value = 70003.75
$ws.Range("L10").Value = 70003.75
$ws.Range("N10").Value = -70343.75

# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 6330.421
$ws.Range("I32").Value = 5831.2446
$ws.Range("J32").Value = 8202.333000000001
$ws.Range("K32").Value = 5831.2446
$ws.Range("L32").Value = 8202.333000000001
$ws.Range("M32").Value = -5544.2446
$ws.Range("N32").Value = -8776.333000000001

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 3187.1404
$ws.Range("I61").Value = 4064.125
$ws.Range("J61").Value = 1123.6471
$ws.Range("K61").Value = 4064.125
$ws.Range("L61").Value = 1123.6471
$ws.Range("M61").Value = -3852.125
$ws.Range("N61").Value = -1547.6471

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 735083.7
$ws.Range("I122").Value = 886600.0600000001
$ws.Range("J122").Value = 2754.6667
$ws.Range("K122").Value = 2659800.18
$ws.Range("L122").Value = 8264.000100000001
$ws.Range("M122").Value = -2657350.18
$ws.Range("N122").Value = -13164.0001

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 4153.2856
$ws.Range("I132").Value = 2060.1428
$ws.Range("J132").Value = 8339.571
$ws.Range("K132").Value = 6180.428400000001
$ws.Range("L132").Value = 25018.713
$ws.Range("M132").Value = -3650.428400000001
$ws.Range("N132").Value = -30078.713

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 3187.1404
$ws.Range("I136").Value = 4064.125
$ws.Range("J136").Value = 1123.6471
$ws.Range("K136").Value = 12192.375
$ws.Range("L136").Value = 3370.9413
$ws.Range("M136").Value = -9642.375
$ws.Range("N136").Value = -8470.941299999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 12 (Leve Item ID 1604)
$ws.Range("H12").Value = 2002.5
$ws.Range("I12").Value = 2002.5
$ws.Range("K12").Value = 2002.5
$ws.Range("M12").Value = -1832.5

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2609.883
$ws.Range("I31").Value = 2230
$ws.Range("J31").Value = 2725.7795
$ws.Range("K31").Value = 2230
$ws.Range("L31").Value = 2725.7795
$ws.Range("M31").Value = -1935
$ws.Range("N31").Value = -3315.7795

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2609.883
$ws.Range("I34").Value = 2230
$ws.Range("J34").Value = 2725.7795
$ws.Range("K34").Value = 2230
$ws.Range("L34").Value = 2725.7795
$ws.Range("M34").Value = -2028
$ws.Range("N34").Value = -3129.7795

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2295.6
$ws.Range("I132").Value = 1365.0476
$ws.Range("J132").Value = 4466.8887
$ws.Range("K132").Value = 4095.142800000001
$ws.Range("L132").Value = 13400.6661
$ws.Range("M132").Value = -1565.142800000001
$ws.Range("N132").Value = -18460.6661

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 12 (Leve Item ID 4854)
$ws.Range("H12").Value = 5882480
$ws.Range("I12").Value = 9091003
$ws.Range("J12").Value = 187.66667
$ws.Range("K12").Value = 27273009
$ws.Range("L12").Value = 563.00001
$ws.Range("M12").Value = -27272836
$ws.Range("N12").Value = -909.00001

# Row 59 (Leve Item ID 4694)
$ws.Range("H59").Value = 6500
$ws.Range("J59").Value = 6500
$ws.Range("L59").Value = 19500
$ws.Range("N59").Value = -20580

# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 978.07574
$ws.Range("I107").Value = 317.625
$ws.Range("J107").Value = 1189.42
$ws.Range("K107").Value = 952.875
$ws.Range("L107").Value = 3568.26
$ws.Range("M107").Value = 967.125
$ws.Range("N107").Value = -7408.26

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 46 (Leve Item ID 2078)
$ws.Range("H46").Value = 7391.375
$ws.Range("J46").Value = 6548.3335
$ws.Range("L46").Value = 6548.3335
$ws.Range("N46").Value = -6860.3335

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 7880.8
$ws.Range("I7").Value = 2600.5833
$ws.Range("J7").Value = 29001.666
$ws.Range("K7").Value = 2600.5833
$ws.Range("L7").Value = 29001.666
$ws.Range("M7").Value = -2488.5833
$ws.Range("N7").Value = -29225.666

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 2913752.2
$ws.Range("I122").Value = 3763772.5
$ws.Range("K122").Value = 11291317.5
$ws.Range("M122").Value = -11288867.5

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 7880.8
$ws.Range("I126").Value = 2600.5833
$ws.Range("J126").Value = 29001.666
$ws.Range("K126").Value = 7801.749899999999
$ws.Range("L126").Value = 87004.99800000001
$ws.Range("M126").Value = -5331.749899999999
$ws.Range("N126").Value = -91944.99800000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 7 (Leve Item ID 2661)
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

# Row 12 (Leve Item ID 3316)
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

# Row 39 (Leve Item ID 3106)
$ws.Range("H39").Value = 18000
$ws.Range("J39").Value = 18000
$ws.Range("L39").Value = 18000
$ws.Range("N39").Value = -18826

# Row 42 (Leve Item ID 3372)
$ws.Range("H42").Value = 23500
$ws.Range("J42").Value = 23500
$ws.Range("L42").Value = 23500
$ws.Range("N42").Value = -24256

# Row 123 (Leve Item ID 34127)
$ws.Range("H123").Value = 35429
$ws.Range("J123").Value = 35429
$ws.Range("L123").Value = 35429
$ws.Range("N123").Value = -45229

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 965.73914
$ws.Range("I126").Value = 746.93335
$ws.Range("K126").Value = 2240.80005
$ws.Range("M126").Value = 229.1999500000002
